# Update Price (column D) and Volume(1h) (column E) figures for the
# symbol list refresh performed by the scheduled GitHub Actions job.
#
# Source cells are stored as literal text (t="inlineStr" in the OOXML),
# not numbers/percentages, so each write forces the cell to Text via
# NumberFormat "@" before assigning the value, then restores the cell's
# style index to the sheet default ("Normal") so no stray number-format
# style is left attached to the cell -- matching the original workbook,
# where none of these data cells carry an `s` (style) attribute.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = "316.47" },
    @{ Row = 2; Col = 5; Value = "4.43%" },
    @{ Row = 3; Col = 4; Value = "48.20" },
    @{ Row = 3; Col = 5; Value = "11.24%" },
    @{ Row = 4; Col = 4; Value = "5.285" },
    @{ Row = 4; Col = 5; Value = "4.61%" },
    @{ Row = 5; Col = 4; Value = "0.07934" },
    @{ Row = 5; Col = 5; Value = "3.32%" },
    @{ Row = 6; Col = 4; Value = "4.600" },
    @{ Row = 6; Col = 5; Value = "4.37%" },
    @{ Row = 7; Col = 4; Value = "1.341" },
    @{ Row = 7; Col = 5; Value = "33.36%" },
    @{ Row = 8; Col = 4; Value = "1.635" },
    @{ Row = 8; Col = 5; Value = "1.84%" },
    @{ Row = 9; Col = 4; Value = "0.1276" },
    @{ Row = 9; Col = 5; Value = "2.43%" },
    @{ Row = 10; Col = 4; Value = "0.1949" },
    @{ Row = 10; Col = 5; Value = "5.05%" },
    @{ Row = 11; Col = 4; Value = "0.09402" },
    @{ Row = 11; Col = 5; Value = "3.14%" },
    @{ Row = 12; Col = 4; Value = "0.04626" },
    @{ Row = 12; Col = 5; Value = "10.92%" },
    @{ Row = 13; Col = 4; Value = "0.1047" },
    @{ Row = 13; Col = 5; Value = "-0.02%" },
    @{ Row = 14; Col = 4; Value = "0.001324" },
    @{ Row = 14; Col = 5; Value = "2.41%" },
    @{ Row = 15; Col = 4; Value = "0.04169" },
    @{ Row = 15; Col = 5; Value = "0.11%" },
    @{ Row = 16; Col = 4; Value = "0.005839" },
    @{ Row = 16; Col = 5; Value = "1.39%" },
    @{ Row = 17; Col = 4; Value = "3.327" },
    @{ Row = 17; Col = 5; Value = "-0.09%" },
    @{ Row = 18; Col = 4; Value = "2.427" },
    @{ Row = 18; Col = 5; Value = "3.02%" },
    @{ Row = 19; Col = 4; Value = "0.3506" },
    @{ Row = 19; Col = 5; Value = "4.65%" },
    @{ Row = 20; Col = 4; Value = "8.070" },
    @{ Row = 20; Col = 5; Value = "-4.42%" },
    @{ Row = 21; Col = 4; Value = "0.1394" },
    @{ Row = 21; Col = 5; Value = "-0.22%" },
    @{ Row = 23; Col = 5; Value = "2.98%" },
    @{ Row = 24; Col = 4; Value = "0.004187" },
    @{ Row = 24; Col = 5; Value = "-6.29%" },
    @{ Row = 25; Col = 4; Value = "0.0001354" },
    @{ Row = 25; Col = 5; Value = "0.35%" },
    @{ Row = 26; Col = 4; Value = "0.0003550" },
    @{ Row = 26; Col = 5; Value = "-95.22%" },
    @{ Row = 38; Col = 4; Value = "0.02656" },
    @{ Row = 38; Col = 5; Value = "8.21%" },
    @{ Row = 39; Col = 4; Value = "0.05698" },
    @{ Row = 39; Col = 5; Value = "7.83%" },
    @{ Row = 40; Col = 4; Value = "0.01079" },
    @{ Row = 40; Col = 5; Value = "80.71%" },
    @{ Row = 41; Col = 4; Value = "0.008027" },
    @{ Row = 41; Col = 5; Value = "4.41%" },
    @{ Row = 42; Col = 4; Value = "0.1436" },
    @{ Row = 42; Col = 5; Value = "6.81%" },
    @{ Row = 43; Col = 4; Value = "0.007704" },
    @{ Row = 43; Col = 5; Value = "4.65%" },
    @{ Row = 44; Col = 4; Value = "0.008488" },
    @{ Row = 44; Col = 5; Value = "13.54%" },
    @{ Row = 45; Col = 4; Value = "0.3153" },
    @{ Row = 45; Col = 5; Value = "4.40%" },
    @{ Row = 46; Col = 4; Value = "0.00006809" },
    @{ Row = 46; Col = 5; Value = "1.45%" },
    @{ Row = 47; Col = 5; Value = "0.31%" },
    @{ Row = 48; Col = 4; Value = "0.05495" },
    @{ Row = 48; Col = 5; Value = "34.23%" },
    @{ Row = 49; Col = 4; Value = "0.004013" },
    @{ Row = 49; Col = 5; Value = "-4.48%" },
    @{ Row = 50; Col = 4; Value = "0.00002106" },
    @{ Row = 50; Col = 5; Value = "0.31%" },
    @{ Row = 51; Col = 5; Value = "0.31%" }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
